# Apply "new covid case, new column" edit to samples-details workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New row 55: a new covid-passport phishing sample (fill most of it first,
#    the two new-URL cells and the facebook link get filled in further down
#    to keep the shared-string table append order matching the source edit).
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "msg"
$ws.Range("C55").Value = "2021-07-31"
$ws.Range("D55").Value = "Official"
$ws.Range("E55").Value = "shortened"
$ws.Range("F55").Value = "opportunity"
$ws.Range("G55").Value = "mt"
$ws.Range("H55").Value = "no"
$ws.Range("I55").Value = "buy covid passport"
$ws.Range("J55").Value = "Government"

# 2. Move the "redirects to <url>" values out of column K into column L as bare URLs,
#    for the rows that currently hold them.
$ws.Range("K46").ClearContents()
$ws.Range("L46").Value = "https://doctorbrew.pl/wp-admin/user/-/"

$ws.Range("K48").ClearContents()
$ws.Range("L48").Value = "https://www.restaurant-apron.at/wp-admin/network/-/"

$ws.Range("K50").ClearContents()
$ws.Range("L50").Value = "https://officehotmail2021.weebly.com/"

$ws.Range("K51").ClearContents()
$ws.Range("L51").Value = "https://watson.pe/onlineBov/"

# 3. Finish off row 55 with its own link + moreinfo values.
$ws.Range("L55").Value = "https://dhl-mt-cliint-srvscs-soynius-rfiid-delivery-sophiechappot211321.codeanyapp.com/mtxx/index1.php"
$ws.Range("K55").Value = "https://www.facebook.com/story.php?story_fbid=4377523138997877&id=488536334563263"

# 4. New column L header: "redirectsto"
$ws.Range("L1").Value = "redirectsto"

# 5. Freeze panes: freeze header row + first column, selection on A2.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2").Select()
